$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24:B30").Style = "Normal"

$ws.Range("B24").Value = "pavan kumar"
$ws.Range("B25").Value = "Narayan Ghorpade"
$ws.Range("B26").Value = "sriharsha bollu"
$ws.Range("B27").Value = "Gajjala Hemanth Kumar Reddy"
$ws.Range("B28").Value = "Prapul Reddy Patlolla"
$ws.Range("B29").Value = "Divya M"
$ws.Range("B30").Value = "Thonduru Arun Kumar"

$ws.Range("B5:B30").Select()
